# Update column G ("K") values for rows 2-40 (data rows) as per the regenerated
# save_data: use K instead of Strike#, after recalculating std/mean and s_vals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(2,2,0,2,0,5,5,7,5,5,10,5,1,6,7,4,4,7,7,6,7,8,4,5,10,10,9,8,3,4,8,7,7,7,5,4,3,2,1)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value = $newValues[$i]
}
